$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure account-number column (C) stays text even though the new values
# are purely numeric strings (matches the source file's t="str" storage).
$ws.Range("C2:C5").NumberFormat = "@"

# --- Row 2 ---
$ws.Range("A2").Value = "SAMIRA TATA"
$ws.Range("B2").Value = "D524564"
$ws.Range("C2").Value = "335463513748543615567464"
$ws.Range("F2").Value = "Supervision"
$ws.Range("G2").Value = "554/SUP FES 1"
$ws.Range("I2").Value = 50000
$ws.Range("J2").Value = 6000
$ws.Range("K2").Value = 44000

# --- Row 3 ---
$ws.Range("A3").Value = "NABIL KAMAL"
$ws.Range("B3").Value = "L3578354"
$ws.Range("C3").Value = "345534544587485743558673"
$ws.Range("D3").Value = "AGG1"
$ws.Range("F3").Value = "Logement de fonction"
$ws.Range("G3").Value = "044/LF/FES VILLE "
$ws.Range("I3").Value = 50000
$ws.Range("J3").Value = 6000
$ws.Range("K3").Value = 44000

# --- Row 4 ---
$ws.Range("A4").Value = "KHADIJA LALA"
$ws.Range("B4").Value = "K5443645"
$ws.Range("C4").Value = "354564564324158786713544"
$ws.Range("D4").Value = "AG 100"
$ws.Range("E4").Value = "BP"
$ws.Range("F4").Value = "Direction régionale"
$ws.Range("G4").Value = "044/FES VILLE "
$ws.Range("I4").Value = 60000
$ws.Range("J4").Value = 9000
$ws.Range("K4").Value = 51000

# --- Row 5 ---
$ws.Range("A5").Value = "KHADIJA LALA"
$ws.Range("B5").Value = "K5443645"
$ws.Range("C5").Value = "354564564324158786713544"
$ws.Range("D5").Value = "AG 100"
$ws.Range("F5").Value = "Direction régionale"
$ws.Range("G5").Value = "044/FES VILLE "
$ws.Range("I5").Value = 20000
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 18500

# --- Row 6 (former row 6 "SAMIA NARA" row becomes the new totals row) ---
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("I6").Value = 180000
$ws.Range("J6").Value = 22500
$ws.Range("K6").Value = 157500

# --- Remove old rows 7-10 (table now ends at row 6) ---
$ws.Range("A7:K10").Delete()
